# Generate Report for Handoff
# Adds a new tracked file (be752219-...) as row 3 on each of the three
# report sheets (Overview, zh-cn, de-de), mirroring the existing
# 16dae7a3-... row, and grows the three ListObjects / dimensions to match.

$wb = $excel.ActiveWorkbook

$mdName    = "be752219-0b85-433d-994b-14e6a030e7e4ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$mdPath    = "e2e\" + $mdName
$mdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e75ddadd9cc92399608d0f03fea06ce2b6e74027/e2e/" + $mdName

$xlfZhCn   = "be752219-0b85-433d-994b-14e6a030e7e4oooooooooooooooooooooooooooooooooooooooo.5f9fd6a092a55807d72e07a7b01ff6579179113b.zh-cn.xlf"
$xlfDeDe   = "be752219-0b85-433d-994b-14e6a030e7e4oooooooooooooooooooooooooooooooooooooooo.5f9fd6a092a55807d72e07a7b01ff6579179113b.de-de.xlf"

$status      = "Ready for handoff"
$genDate     = "2016-08-31 16:31:39"
$handoffZh   = "2016-08-31 16:31:35"
$handoffDe   = "2016-08-31 16:31:39"
$noHandback  = "0001-01-01 00:00:00"
$dtFormat    = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1) -> row 3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(3, 1).Value = $mdName
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(3, 2), $mdUrl, "", "", $mdPath) | Out-Null
$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 5).Value = $status
$wsOverview.Cells.Item(3, 6).Value = $status
$wsOverview.Cells.Item(3, 7).Value = $genDate
$wsOverview.Cells.Item(3, 7).NumberFormat = $dtFormat

$wsOverview.Columns.Item(5).ColumnWidth = 16.26
$wsOverview.Columns.Item(6).ColumnWidth = 16.26

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) -> row 3
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(3, 1), $mdUrl, "", "", $mdName) | Out-Null
$wsZhCn.Cells.Item(3, 2).Value = ".md"
$wsZhCn.Cells.Item(3, 3).Value = $status
$wsZhCn.Cells.Item(3, 4).Value = "e2e"
$wsZhCn.Cells.Item(3, 5).Value = "ht"
$wsZhCn.Cells.Item(3, 6).Value = "'False"
$wsZhCn.Cells.Item(3, 7).Value = $xlfZhCn
$wsZhCn.Cells.Item(3, 8).Value = $handoffZh
$wsZhCn.Cells.Item(3, 8).NumberFormat = $dtFormat
$wsZhCn.Cells.Item(3, 11).Value = $noHandback
$wsZhCn.Cells.Item(3, 11).NumberFormat = $dtFormat
$wsZhCn.Cells.Item(3, 13).Value = "'True"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.26

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3) -> row 3
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(3, 1), $mdUrl, "", "", $mdName) | Out-Null
$wsDeDe.Cells.Item(3, 2).Value = ".md"
$wsDeDe.Cells.Item(3, 3).Value = $status
$wsDeDe.Cells.Item(3, 4).Value = "e2e"
$wsDeDe.Cells.Item(3, 5).Value = "ht"
$wsDeDe.Cells.Item(3, 6).Value = "'False"
$wsDeDe.Cells.Item(3, 7).Value = $xlfDeDe
$wsDeDe.Cells.Item(3, 8).Value = $status
$wsDeDe.Cells.Item(3, 11).Value = $noHandback
$wsDeDe.Cells.Item(3, 11).NumberFormat = $dtFormat
$wsDeDe.Cells.Item(3, 13).Value = "'True"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.26

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
